$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Primera, Metropolitana): date 44301 -> 44322, volumes/prices change
$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("S4").Value = 1714

# Row 5 (Segunda, Metropolitana): date 44301 -> 44322
$ws.Range("D5").Value = 44322
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("S5").Value = 1143

# Row 8 (Primera, Metropolitana): date 44322 -> 44302
$ws.Range("D8").Value = 44302
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 2143

# Row 9 (Segunda, Metropolitana): date 44322 -> 44302
$ws.Range("D9").Value = 44302
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 12 (Primera, Provincia de Santiago -> Región Metropolitana): date 44299 -> 44301
$ws.Range("D12").Value = 44301
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 14000
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 2000

# Row 13 (Segunda, Provincia de Santiago -> Región Metropolitana): date 44299 -> 44301
$ws.Range("D13").Value = 44301
$ws.Range("M13").Value = 80
$ws.Range("R13").Value = "Región Metropolitana"

# Row 14 (Primera, Región Metropolitana -> Provincia de Santiago): date 44302 -> 44299
$ws.Range("D14").Value = 44299
$ws.Range("M14").Value = 80
$ws.Range("R14").Value = "Provincia de Santiago"

# Row 15 (Segunda, Región Metropolitana -> Provincia de Santiago): date 44302 -> 44299
$ws.Range("D15").Value = 44299
$ws.Range("M15").Value = 75
$ws.Range("R15").Value = "Provincia de Santiago"
